$d = $word.ActiveDocument

# Remove the paragraph that holds the "<<time_placeholder>>" field,
# deleting its run content as well as the paragraph mark so the
# remaining "<<date_placeholder>>" paragraph becomes the last one
# in the body (immediately followed by the section properties).
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*time_placeholder*") {
        $p.Range.Delete()
    }
}
